# Push Notification Support + UI/Backend Enhancements
# Visitors sheet data update:
#  - Remove the rejected "Pancholi dn" visitor row.
#  - The previously-approved "Harsh" row moves up to row 2, but loses its
#    "Reason" value (column C is now blank for every data row).
#  - Two brand-new "Harsh" visits are appended (row 3 = pending, row 4 = approved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Visitors")

# Clear out the old data rows and rebuild rows 2:4 from scratch.
$ws.Range("A2:F3").ClearContents()

# Phone numbers must stay text (as in the source data) rather than get
# coerced to numbers. Temporarily force Text format for the write, then
# restore the default "Normal" style so no stray per-cell style sticks.
$ws.Range("B2:B4").NumberFormat = "@"

# Row 2 - existing approved visit (kept, Reason cleared)
$ws.Cells.Item(2, 1).Value = "Harsh"
$ws.Cells.Item(2, 2).Value = "1234567890"
$ws.Cells.Item(2, 4).Value = "approved"
$ws.Cells.Item(2, 5).Value = "4/7/2025, 6:01:11 pm"
$ws.Cells.Item(2, 6).Value = "https://res.cloudinary.com/drdw2abup/image/upload/v1751632271/visitors/gdet8s6acynxjzp5tibe.jpg"

# Row 3 - new pending visit
$ws.Cells.Item(3, 1).Value = "Harsh"
$ws.Cells.Item(3, 2).Value = "2134224231"
$ws.Cells.Item(3, 4).Value = "pending"
$ws.Cells.Item(3, 5).Value = "5/7/2025, 12:37:40 pm"
$ws.Cells.Item(3, 6).Value = "https://res.cloudinary.com/drdw2abup/image/upload/v1751699260/visitors/s8lx4buuoboaabtz39lh.jpg"

# Row 4 - new approved visit
$ws.Cells.Item(4, 1).Value = "Harsh"
$ws.Cells.Item(4, 2).Value = "3458365863"
$ws.Cells.Item(4, 4).Value = "approved"
$ws.Cells.Item(4, 5).Value = "5/7/2025, 12:40:26 pm"
$ws.Cells.Item(4, 6).Value = "https://res.cloudinary.com/drdw2abup/image/upload/v1751699426/visitors/y2lqbutnk4n2xgjvbtmk.jpg"

$ws.Range("B2:B4").Style = "Normal"
